$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 0.225007
$ws.Range("H2").Value = 0.675021
$ws.Range("I2").Value = 0.0376013806128968
$ws.Range("J2").Value = 0.0376013806128968
$ws.Range("M2").Value = 168.1098273333333
$ws.Range("N2").Value = 504.329482
$ws.Range("O2").Value = 0.2984182258032519
$ws.Range("P2").Value = 0.298418225803252
$ws.Range("Q2").Value = 37.82588791879133
$ws.Range("R2").Value = 340.432991269122
$ws.Range("S2").Value = 0.01122093729025346
$ws.Range("T2").Value = 0.01122093729025346
$ws.Range("G3").Value = 0.225007
$ws.Range("H3").Value = 0.675021
$ws.Range("I3").Value = 0.0376013806128968
$ws.Range("J3").Value = 0.0376013806128968
$ws.Range("O3").Value = 0.2893586437755394
$ws.Range("P3").Value = 0.2893586437755394
$ws.Range("Q3").Value = 36.67754406864967
$ws.Range("R3").Value = 330.097896617847
$ws.Range("S3").Value = 0.01088028449823568
$ws.Range("T3").Value = 0.01088028449823568
$ws.Range("G4").Value = 0.225007
$ws.Range("H4").Value = 0.675021
$ws.Range("I4").Value = 0.0376013806128968
$ws.Range("J4").Value = 0.0376013806128968
$ws.Range("M4").Value = 165.99353
$ws.Range("N4").Value = 497.98059
$ws.Range("O4").Value = 0.294661504941043
$ws.Range("P4").Value = 0.294661504941043
$ws.Range("Q4").Value = 37.34970620471
$ws.Range("R4").Value = 336.14735584239
$ws.Range("S4").Value = 0.01107967939925713
$ws.Range("T4").Value = 0.01107967939925713
$ws.Range("G5").Value = 0.225007
$ws.Range("H5").Value = 0.675021
$ws.Range("I5").Value = 0.0376013806128968
$ws.Range("J5").Value = 0.0376013806128968
$ws.Range("M5").Value = 66.22673433333334
$ws.Range("N5").Value = 198.680203
$ws.Range("O5").Value = 0.1175616254801657
$ws.Range("P5").Value = 0.1175616254801657
$ws.Range("Q5").Value = 14.90147881214033
$ws.Range("R5").Value = 134.113309309263
$ws.Range("S5").Value = 0.004420479425150535
$ws.Range("T5").Value = 0.004420479425150535
$ws.Range("G6").Value = 4.404016666666667
$ws.Range("I6").Value = 0.7359642451518149
$ws.Range("J6").Value = 0.7359642451518149
$ws.Range("M6").Value = 168.1098273333333
$ws.Range("N6").Value = 504.329482
$ws.Range("O6").Value = 0.2984182258032519
$ws.Range("P6").Value = 0.298418225803252
$ws.Range("Q6").Value = 740.3584814064557
$ws.Range("R6").Value = 6663.226332658101
$ws.Range("S6").Value = 0.2196251442928342
$ws.Range("T6").Value = 0.2196251442928342
$ws.Range("G7").Value = 4.404016666666667
$ws.Range("I7").Value = 0.7359642451518149
$ws.Range("J7").Value = 0.7359642451518149
$ws.Range("O7").Value = 0.2893586437755394
$ws.Range("P7").Value = 0.2893586437755394
$ws.Range("Q7").Value = 717.8821786465945
$ws.Range("R7").Value = 6460.93960781935
$ws.Range("S7").Value = 0.2129576158444178
$ws.Range("T7").Value = 0.2129576158444178
$ws.Range("G8").Value = 4.404016666666667
$ws.Range("I8").Value = 0.7359642451518149
$ws.Range("J8").Value = 0.7359642451518149
$ws.Range("M8").Value = 165.99353
$ws.Range("N8").Value = 497.98059
$ws.Range("O8").Value = 0.294661504941043
$ws.Range("P8").Value = 0.294661504941043
$ws.Range("Q8").Value = 731.0382726788333
$ws.Range("R8").Value = 6579.3444541095
$ws.Range("S8").Value = 0.2168603320592325
$ws.Range("T8").Value = 0.2168603320592325
$ws.Range("G9").Value = 4.404016666666667
$ws.Range("I9").Value = 0.7359642451518149
$ws.Range("J9").Value = 0.7359642451518149
$ws.Range("M9").Value = 66.22673433333334
$ws.Range("N9").Value = 198.680203
$ws.Range("O9").Value = 0.1175616254801657
$ws.Range("P9").Value = 0.1175616254801657
$ws.Range("Q9").Value = 291.6636417829056
$ws.Range("R9").Value = 2624.972776046151
$ws.Range("S9").Value = 0.08652115295533049
$ws.Range("T9").Value = 0.08652115295533049
$ws.Range("I10").Value = 0.04918519786094349
$ws.Range("J10").Value = 0.0491851978609435
$ws.Range("M10").Value = 168.1098273333333
$ws.Range("N10").Value = 504.329482
$ws.Range("O10").Value = 0.2984182258032519
$ws.Range("P10").Value = 0.298418225803252
$ws.Range("Q10").Value = 49.47886889327421
$ws.Range("R10").Value = 445.3098200394679
$ws.Range("S10").Value = 0.01467775948144466
$ws.Range("T10").Value = 0.01467775948144466
$ws.Range("I11").Value = 0.04918519786094349
$ws.Range("J11").Value = 0.0491851978609435
$ws.Range("O11").Value = 0.2893586437755394
$ws.Range("P11").Value = 0.2893586437755394
$ws.Range("S11").Value = 0.01423216214687417
$ws.Range("T11").Value = 0.01423216214687417
$ws.Range("I12").Value = 0.04918519786094349
$ws.Range("J12").Value = 0.0491851978609435
$ws.Range("M12").Value = 165.99353
$ws.Range("N12").Value = 497.98059
$ws.Range("O12").Value = 0.294661504941043
$ws.Range("P12").Value = 0.294661504941043
$ws.Range("Q12").Value = 48.85599038607332
$ws.Range("R12").Value = 439.7039134746599
$ws.Range("S12").Value = 0.01449298442252858
$ws.Range("T12").Value = 0.01449298442252858
$ws.Range("I13").Value = 0.04918519786094349
$ws.Range("J13").Value = 0.0491851978609435
$ws.Range("M13").Value = 66.22673433333334
$ws.Range("N13").Value = 198.680203
$ws.Range("O13").Value = 0.1175616254801657
$ws.Range("P13").Value = 0.1175616254801657
$ws.Range("Q13").Value = 19.49216150708022
$ws.Range("R13").Value = 175.429453563722
$ws.Range("S13").Value = 0.005782291810096084
$ws.Range("T13").Value = 0.005782291810096084
$ws.Range("G14").Value = 1.060660666666667
$ws.Range("H14").Value = 3.181982
$ws.Range("I14").Value = 0.1772491763743448
$ws.Range("J14").Value = 0.1772491763743448
$ws.Range("M14").Value = 168.1098273333333
$ws.Range("N14").Value = 504.329482
$ws.Range("O14").Value = 0.2984182258032519
$ws.Range("P14").Value = 0.298418225803252
$ws.Range("Q14").Value = 178.3074815325916
$ws.Range("R14").Value = 1604.767333793324
$ws.Range("S14").Value = 0.05289438473871965
$ws.Range("T14").Value = 0.05289438473871966
$ws.Range("G15").Value = 1.060660666666667
$ws.Range("H15").Value = 3.181982
$ws.Range("I15").Value = 0.1772491763743448
$ws.Range("J15").Value = 0.1772491763743448
$ws.Range("O15").Value = 0.2893586437755394
$ws.Range("P15").Value = 0.2893586437755394
$ws.Range("Q15").Value = 172.8943025930305
$ws.Range("R15").Value = 1556.048723337274
$ws.Range("S15").Value = 0.05128858128601179
$ws.Range("T15").Value = 0.05128858128601179
$ws.Range("G16").Value = 1.060660666666667
$ws.Range("H16").Value = 3.181982
$ws.Range("I16").Value = 0.1772491763743448
$ws.Range("J16").Value = 0.1772491763743448
$ws.Range("M16").Value = 165.99353
$ws.Range("N16").Value = 497.98059
$ws.Range("O16").Value = 0.294661504941043
$ws.Range("P16").Value = 0.294661504941043
$ws.Range("Q16").Value = 176.0628081921533
$ws.Range("R16").Value = 1584.56527372938
$ws.Range("S16").Value = 0.0522285090600248
$ws.Range("T16").Value = 0.0522285090600248
$ws.Range("G17").Value = 1.060660666666667
$ws.Range("H17").Value = 3.181982
$ws.Range("I17").Value = 0.1772491763743448
$ws.Range("J17").Value = 0.1772491763743448
$ws.Range("M17").Value = 66.22673433333334
$ws.Range("N17").Value = 198.680203
$ws.Range("O17").Value = 0.1175616254801657
$ws.Range("P17").Value = 0.1175616254801657
$ws.Range("Q17").Value = 70.24409218914957
$ws.Range("R17").Value = 632.196829702346
$ws.Range("S17").Value = 0.02083770128958855
$ws.Range("T17").Value = 0.02083770128958855
